$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "THÀNH VIÊN" (member) column for the design sub-tasks under
# "2.2 Thiết kế giao diện" (rows 9-12 -> Đăng) and "2.3 Thiết kế dữ liệu"
# (rows 13-16 -> Đạt).
foreach ($r in 9..12) {
    $ws.Range("E$r").Value = "Lê Nguyễn Hoài Đăng"
}
foreach ($r in 13..16) {
    $ws.Range("E$r").Value = "Lê Phước Anh Đạt"
}

# Update the saved view: scroll the window so row 7 is at the top and the
# current selection is B10 (matches the author's working position while
# filling in the plan).
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B10").Select()
